$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns: one before old column C, one before old column D
# (which, after the first insert, sits at column E).
$ws.Columns("C").Insert()
$ws.Columns("E").Insert()

# The new columns repeat the header text of the column to their left.
$ws.Range("C1").Value = $ws.Range("B1").Value2
$ws.Range("E1").Value = $ws.Range("D1").Value2

# New column C: = old B - 100 (growth over 1900 baseline)
$ws.Range("C2").Formula = "=B2-100"
$ws.Range("C3:C11").Formula = "=B3-100"

# New column E: = old D - 100 (growth over 1900 baseline)
$ws.Range("E2").Formula = "=D2-100"
$ws.Range("E3:E11").Formula = "=D3-100"

# New data (insects and snails habitat shares) for rows 3-7, columns G:I.
# These are new historical estimates, so mark them with the red font style.
$newData = @(
    @(3, 96, 93, 90),
    @(4, 90, 80, 75),
    @(5, 85, 65, 60),
    @(6, 80, 45, 40),
    @(7, 70, 30, 20)
)

foreach ($entry in $newData) {
    $r = $entry[0]
    $ws.Cells.Item($r, 7).Value = $entry[1]
    $ws.Cells.Item($r, 8).Value = $entry[2]
    $ws.Cells.Item($r, 9).Value = $entry[3]
    $ws.Range("G" + $r + ":I" + $r).Font.Color = 255
}

# Row 2 already has values 100 for columns G:I by virtue of the column shift,
# but H2/I2 did not exist before (old sheet stopped at G) - add them explicitly.
$ws.Range("H2").Value = 100
$ws.Range("I2").Value = 100

# New summary row with a simple formula.
$ws.Range("I17").Formula = "=94/5"

# Update the selection to match the saved workbook view.
$ws.Range("A2:A11").Select()
